$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: D2 becomes a TRUE() formula instead of a hard-coded boolean ---
$ws.Range("D2").Formula = "=TRUE()"

# --- Row 3: new history row, mirroring row 2 but with the new service account ---
$ws.Range("A3").Value = $ws.Range("A2").Value2
$ws.Range("B3").Value = $ws.Range("B2").Value2
$ws.Range("C3").Value = "Service-account-mosip-resident-client"

$ws.Range("D3").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("D3").Formula = "=TRUE()"

$ws.Range("E3").Value = $ws.Range("E2").Value2

# --- Selection moves to E3 ---
$ws.Range("E3").Select() | Out-Null
